$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "MSc AI/DS 입학시험 후기"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msc-ai-ds-prep-review/#utm_source=rss&utm_medium=rss&utm_campaign=msc-ai-ds-prep-review"

$ws.Range("D28").Value = "[논문리뷰] Momentum Observer-Based Collision Detection Using LSTM for Model Uncertainty Learning"
$ws.Range("E28").Value = "https://ropiens.tistory.com/186"

$ws.Range("D46").Value = "매독(syphilis) 임상양상"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/460"

$ws.Range("D51").Value = "[python+pandas] 판다스 데이터프레임 loc, at, iloc, iat 메소드 비교"
$ws.Range("E51").Value = "https://bskyvision.com/1221"
